$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "shyam ghosh" record (row 3) replaces the "Alice1 Johnson" record
# (row 2). Copy the whole row across first so every cell keeps its original
# type/formatting (this also keeps customer_id/phone_number/country_code as
# text instead of Excel re-interpreting the numeric-looking strings as
# numbers).
$ws.Range("A3:O3").Copy($ws.Range("A2:O2"))

# The call_id on the surviving row is refreshed to 178.
$ws.Range("A2").Value = 178

# Notes / tasks (columns M, N) are removed from the LLM prompt context.
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""

# The old row 3 (now folded into row 2) is no longer needed.
$ws.Rows.Item(3).Delete()
